$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New "Save" header column (H), matching the formatting used by the
# existing header cells (e.g. G1: bold, bordered, centered).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Fill the new column's data rows (H2:H8) with 0, matching the plain
# (unstyled) numeric cells used elsewhere in the data rows.
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}
